$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '50.835.12'
$ws.Range('E2').Value = '  -0.85%  '

$ws.Range('D3').Value = '2.927.32'
$ws.Range('E3').Value = '  -1.48%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '375.38'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.64%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '100.16'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.61%  '

$ws.Range('E7').Value = '  -0.83%  '

$ws.Range('E8').Value = '  -0.07%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.41%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.00'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.10%  '

$ws.Range('E11').Value = '  -0.52%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0845'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.63%  '

$ws.Range('D13').Value = '3.388.05'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17.96'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.81%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.59'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.15%  '

$ws.Range('E16').Value = '  +53.88%  '

$ws.Range('D17').Value = '2.921.54'
$ws.Range('E17').Value = '  -1.70%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.990'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.06%  '

$ws.Range('D19').Value = '50.786.11'
$ws.Range('E19').Value = '  -0.80%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.06'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -5.72%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.37'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.82%  '

$ws.Range('D22').Value = '0.0₃0954'
$ws.Range('E22').Value = '  -0.04%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.84'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '264.96'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.74%  '

$ws.Range('E25').Value = '  +9.58%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.01'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.60%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.42'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.40%  '

$ws.Range('E28').Value = '  +0.00%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '25.50'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.22%  '

$ws.Range('E30').Value = '  -5.83%  '

$ws.Range('E31').Value = '  -3.45%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '9.98'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.84%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '50.84'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.12%  '

$ws.Range('E34').Value = '  -0.65%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '33.08'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.31%  '

$ws.Range('E36').Value = '  -3.16%  '

$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.09'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +4.35%  '

$ws.Range('E39').Value = '  -0.22%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.44'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.91%  '

$ws.Range('E41').Value = '  +0.00%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.46'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.63%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '119.38'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.65%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.00'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.27%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.38'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.25%  '

$ws.Range('E46').Value = '  -1.94%  '

$ws.Range('E47').Value = '  -1.32%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.267'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.42%  '

$ws.Range('D49').Value = '1.987.84'
$ws.Range('E49').Value = '  -1.90%  '

$ws.Range('E50').Value = '  -2.41%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.19'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.03%  '
